# Auto-generated edit script applying the cryptos.xlsx price/volume update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'51.791.59"
$ws.Range('E2').Value = '  +0.13%  '

$ws.Range('D3').Value = "'2.803.97"
$ws.Range('E3').Value = '  +0.86%  '

$ws.Range('E4').Value = '  +0.02%  '

$ws.Range('D5').Value = "'355.01"
$ws.Range('E5').Value = '  -0.34%  '

$ws.Range('D6').Value = "'111.36"
$ws.Range('E6').Value = '  +1.78%  '

$ws.Range('D7').Value = "'0.559"
$ws.Range('E7').Value = '  +1.02%  '

$ws.Range('E8').Value = '  +0.05%  '

$ws.Range('E9').Value = '  +7.94%  '

$ws.Range('D10').Value = "'40.42"
$ws.Range('E10').Value = '  +2.04%  '

$ws.Range('D12').Value = "'0.0840"
$ws.Range('E12').Value = '  -0.38%  '

$ws.Range('D13').Value = "'20.05"

$ws.Range('E14').Value = '  +2.79%  '

$ws.Range('D15').Value = "'3.249.86"
$ws.Range('E15').Value = '  +1.06%  '

$ws.Range('D16').Value = "'2.811.02"
$ws.Range('E16').Value = '  +0.76%  '

$ws.Range('E17').Value = '  +1.74%  '

$ws.Range('D18').Value = "'51.797.63"
$ws.Range('E18').Value = '  +0.23%  '

$ws.Range('D19').Value = "'7.66"
$ws.Range('E19').Value = '  +1.93%  '

$ws.Range('D20').Value = "'3.19"
$ws.Range('E20').Value = '  +2.79%  '

$ws.Range('D21').Value = "'13.64"
$ws.Range('E21').Value = '  +3.51%  '

$ws.Range('E22').Value = '  +0.99%  '

$ws.Range('D23').Value = "'70.55"
$ws.Range('E23').Value = '  +0.61%  '

$ws.Range('D24').Value = "'268.96"
$ws.Range('E24').Value = '  +0.40%  '

$ws.Range('E25').Value = '  +1.63%  '

$ws.Range('E26').Value = '  +0.09%  '

$ws.Range('E27').Value = '  -0.65%  '

$ws.Range('E28').Value = '  -2.11%  '

$ws.Range('D29').Value = "'38.77"
$ws.Range('E29').Value = '  +11.54%  '

$ws.Range('D30').Value = "'10.40"
$ws.Range('E30').Value = '  +2.13%  '

$ws.Range('D31').Value = "'2.25"
$ws.Range('E31').Value = '  +3.97%  '

$ws.Range('D32').Value = "'52.42"
$ws.Range('E32').Value = '  +1.44%  '

$ws.Range('E33').Value = '  +0.48%  '

$ws.Range('E34').Value = '  +8.95%  '

$ws.Range('D35').Value = "'0.0888"
$ws.Range('E35').Value = '  +6.36%  '

$ws.Range('E36').Value = '  -0.76%  '

$ws.Range('D37').Value = "'0.999"
$ws.Range('E37').Value = '  -0.05%  '

$ws.Range('D38').Value = "'18.80"
$ws.Range('E38').Value = '  -0.25%  '

$ws.Range('E39').Value = '  +2.41%  '

$ws.Range('E40').Value = '  +0.72%  '

$ws.Range('E41').Value = '  +1.32%  '

$ws.Range('E42').Value = '  -2.15%  '

$ws.Range('D43').Value = "'120.46"
$ws.Range('E43').Value = '  +0.62%  '

$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D44').Value = "'22.11"
$ws.Range('E44').Value = '  +1.96%  '

$ws.Range('B45').Value = 'WEMIXToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D45').Value = "'2.21"
$ws.Range('E45').Value = '  -0.05%  '

$ws.Range('D46').Value = "'3.42"
$ws.Range('E46').Value = '  +5.15%  '

$ws.Range('D47').Value = "'2.113.09"
$ws.Range('E47').Value = '  +0.99%  '

$ws.Range('D48').Value = "'2.42"
$ws.Range('E48').Value = '  +6.64%  '

$ws.Range('D49').Value = "'0.951"
$ws.Range('E49').Value = '  +0.31%  '

$ws.Range('E50').Value = '  -1.00%  '

$ws.Range('E51').Value = '  +6.85%  '
